# SprintHours.xlsx update
# - Fill in previously-blank "Hours Spent" (column E) values for the
#   "15/2/17 - 22/2/17" sprint block (rows 268-292)
# - Append two new task blocks / a new sprint header block for
#   "15/3/17 - 22/3/17" (rows 294-319)
# - Move the active selection to reflect where the user ended up working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Hours($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

# ---------------------------------------------------------------------
# Existing rows: fill in the Hours Spent (column E) values that were
# previously left blank.
# ---------------------------------------------------------------------
Set-Hours "E273" 0
Set-Hours "E274" 1
Set-Hours "E275" 1.5

Set-Hours "E277" 4
Set-Hours "E278" 0.5
Set-Hours "E279" 0.5
Set-Hours "E280" 0.5
Set-Hours "E281" 1
Set-Hours "E282" 0.5
Set-Hours "E283" 0.5
Set-Hours "E284" 0.5
Set-Hours "E285" 1
Set-Hours "E286" 1.5

Set-Hours "E288" 0
Set-Hours "E289" 2
Set-Hours "E290" 2
Set-Hours "E291" 1
Set-Hours "E292" 1.5

# ---------------------------------------------------------------------
# New rows 294-319: a fresh "Group Member / Task / Hours Spent / Sprint /
# Total Hours" header plus a new sprint ("15/3/17 - 22/3/17") and its
# per-person task breakdown.
# ---------------------------------------------------------------------

# Header row (same layout as row 267)
$ws.Range("C294").Value = "Group Member "
$ws.Range("E294").Value = "Hours Spent"
$ws.Range("F294").Value = "Sprint "
$ws.Range("G294").Value = "Total Hours"

# Sprint-total row
$ws.Range("C295").Value = "15/3/17 - 22/3/17"
$ws.Range("F295").Value = 20
$ws.Range("G295").Formula = "=E296+E297+E298+E299+E300+E301+E302"

# Matthew Allum's tasks
$ws.Range("C296").Value = "Matthew Allum"
$ws.Range("D296").Value = "As a designer, fix the missing art on particular furniture models"
$ws.Range("D297").Value = "As a designer, make the lighting for the game dynamic (see Github comments)"
$ws.Range("D298").Value = "As a designer, edit the coin values to allow the player more time when cranked in the fuse box"
$ws.Range("D299").Value = "As a designer, make the rooms lighter when the lights go out"
$ws.Range("D300").Value = "As a modeler, retexture the fuse box"
$ws.Range("D301").Value = "As a modeler, model and texture some temporary blinds for the windows"
$ws.Range("D302").Value = "As a group, make a blog post"

# Aaron Mulligan's tasks
$ws.Range("C304").Value = "Aaron Mulligan "
$ws.Range("D304").Value = "As a designer, fix the text prompt when you go up to the fuse box with a coin"
$ws.Range("E304").Value = 0.5
$ws.Range("G304").Formula = "=E313+E307+E308+E309+E310+E311+E312+E306+E305+E304"

$ws.Range("D305").Value = "As a designer, redesign the noticeboard (see Github comments)"
$ws.Range("D306").Value = "As a designer, replace memory 1 image with something more fitting"
$ws.Range("D307").Value = "As a sound artist, create or find music piece for the music player to play"
$ws.Range("D308").Value = "As a sound artist, create sound for when players fall down the stairs"
$ws.Range("D309").Value = "As a designer, make the coin particles stand out more"

$ws.Range("D310").Value = "As a coder, fix the carer model animation warpng backwards when the animation finishes"
$ws.Range("E310").Value = 0.5

$ws.Range("D311").Value = "As a group, make a blog post"
$ws.Range("E311").Value = 1

# Lee Hatchman's tasks
$ws.Range("C315").Value = "Lee Hatchman"
$ws.Range("D315").Value = "As a modeler, model and texture ceiling lights"
$ws.Range("G315").Formula = "=E315+E316+E317+E318+E319"

$ws.Range("D316").Value = "As a modeler, model and texture doors"
$ws.Range("D317").Value = "As a group, make a blog post"

# ---------------------------------------------------------------------
# Apply the "Check Cell" style (the double-border / bold-white-on-grey
# look used throughout this table) and row heights to every new row so
# they visually match the rest of the sheet.
# ---------------------------------------------------------------------
$newRowsRange = $ws.Range("C294:G318")
$newRowsRange.Style = "Check Cell"

for ($r = 294; $r -le 318; $r++) {
    $ws.Rows.Item($r).RowHeight = 16.5
}
$ws.Rows.Item(319).RowHeight = 15.75

# ---------------------------------------------------------------------
# Move the selection to where the editor ended up (bottom of sheet).
# ---------------------------------------------------------------------
$ws.Range("H308").Select()
